# minimal test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Core Java 8")
$ws.Range("A460").Value = "test"
